$wb = $excel.ActiveWorkbook
$nbsp = [char]0x00A0

# ---------------------------------------------------------------------------
# Engine sheet (sheet1): give Q2 a bottom border, move selection, deactivate
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Engine")
$ws1.Range("Q2").Borders(9).Weight = 2
$ws1.Range("H13").Select()

# ---------------------------------------------------------------------------
# Paint sheet (sheet3): re-number the Id column from 0, add a small
# "Categories" lookup table in column G, widen column G, move selection
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Paint")

$ws3.Range("B3").Value = 0
$ws3.Range("B4").Value = 1
$ws3.Range("B5").Value = 2
$ws3.Range("B6").Value = 3
$ws3.Range("B7").Value = 4
$ws3.Range("B8").Value = 5
$ws3.Range("B9").Value = 6
$ws3.Range("B10").Value = 7
$ws3.Range("B11").Value = 8

$ws3.Columns.Item(7).ColumnWidth = 19.6

$ws3.Range("G15").Value = "Categories"
$ws3.Range("G15").Font.Bold = $true
$ws3.Range("G16").Value = "Uni"
$ws3.Range("G17").Value = "Metalic"
$ws3.Range("G18").Value = "Pearlescent" + $nbsp + "paint"

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

$ws3.Range("B12").Select()

# ---------------------------------------------------------------------------
# New Rims sheet (sheet4): id / size / price table
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws4.Name = "Rims"

$ws4.Range("B2").Value = "ID"
$ws4.Range("B2").Font.Bold = $true
$ws4.Range("C2").Value = "Size"
$ws4.Range("C2").Font.Bold = $true
$ws4.Range("D2").Value = "Price"
$ws4.Range("D2").Font.Bold = $true
$ws4.Range("E2").Font.Bold = $true

$ws4.Range("B3").Value = 0
$ws4.Range("C3").Value = 16
$ws4.Range("D3").Value = 0

$ws4.Range("B4").Value = 1
$ws4.Range("C4").Value = 17
$ws4.Range("D4").Value = 600

$ws4.Range("B5").Value = 2
$ws4.Range("C5").Value = 18
$ws4.Range("D5").Value = 1400

$ws4.Range("B6").Value = 3
$ws4.Range("C6").Value = 19
$ws4.Range("D6").Value = 2400

$ws4.PageSetup.TopMargin = 56.692913399999995
$ws4.PageSetup.BottomMargin = 56.692913399999995

$ws4.Activate()
$ws4.Range("D7").Select()
